$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Associated Lines" (column S) values for rows 9-13
$ws.Range("S9").Value  = "Island Line & Kwun Tong Line"
$ws.Range("S10").Value = "South Island Line, Island Line, Walkable, Kwun Tong Line & West Rail Line"
$ws.Range("S11").Value = "Airport Express, Walkable & Disneyland Resort Line"
$ws.Range("S12").Value = "Tsuen Wan Line & East Rail Line"
$ws.Range("S13").Value = "Airport Express, Tsuen Wan Line, Island Line & Tung Chung Line"

# Scroll the frozen (right) pane so its top-left visible cell is B1 instead of O1,
# and move the active selection to E26 (last row the user was entering data into)
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("E26").Select()
